# AINT512 Project Proposal - v1.4.3 edits
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "1.2 - Speech To Text Design" heading: merge the split "To" run back in
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("1.2 – Speech To Text Design", $true, $false, $false, $false, $false, $true, 1, $false, "1.2 – Speech To Text Design", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Grammar paragraph: merge "; an " back into the preceding sentence
# ---------------------------------------------------------------------------
$semi = "several minutes" + [char]8217 + " worth of processing; an "
$d.Content.Find.Execute($semi, $true, $false, $false, $false, $false, $true, 1, $false, $semi, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "2.1 - Basics:" heading: merge "2.1" + " - " together
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2.1 – ", $true, $false, $false, $false, $false, $true, 1, $false, "2.1 – ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Dialogue Specification:" paragraph: merge the lone space run into the
#    following sentence
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("The dialogue specification, or narrative paths, can be seen in ", $true, $false, $false, $false, $false, $true, 1, $false, "The dialogue specification, or narrative paths, can be seen in ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "figure x." : remove proofErr wrap around the full stop
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("figure x.", $true, $false, $false, $false, $false, $true, 1, $false, "figure x.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Remove the old `_GoBack` bookmark sitting before "Personality Specification"
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 7) "Personality Specification:" merge the trailing colon run in
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Personality Specification:", $true, $false, $false, $false, $false, $true, 1, $false, "Personality Specification:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) "2.2 - Grammar Design:" heading: drop the proofErr wrap around "2.2"
#    ("2." and "2" stay separate runs, only the proofErr markers go away)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2.2 – Grammar Design:", $true, $false, $false, $false, $false, $true, 1, $false, "2.2 – Grammar Design:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 9) Grammar Design paragraph: merge the proofErr-wrapped "context" back in
# ---------------------------------------------------------------------------
$ctx = "This will provide enough context to infer the user" + [char]8217 + "s selection."
$d.Content.Find.Execute($ctx, $true, $false, $false, $false, $false, $true, 1, $false, $ctx, 2) | Out-Null

# ---------------------------------------------------------------------------
# 10-13) Drop proofErr wraps around the remaining numbered sub-heading runs
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2.3 – ", $true, $false, $false, $false, $false, $true, 1, $false, "2.3 – ", 2) | Out-Null
$d.Content.Find.Execute("2.4 – ", $true, $false, $false, $false, $false, $true, 1, $false, "2.4 – ", 2) | Out-Null
$d.Content.Find.Execute("2.5 – ", $true, $false, $false, $false, $false, $true, 1, $false, "2.5 – ", 2) | Out-Null
$d.Content.Find.Execute("2.6 – ", $true, $false, $false, $false, $false, $true, 1, $false, "2.6 – ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 14) Final paragraph - real content edit: rewrite the "director" example and
#     relocate the `_GoBack` bookmark to sit right after the new phrase.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("explicitly state the character by name or by pronoun", $true, $false, $false, $false, $false, $true, 1, $false, "explicitly instructs the character to do something", 2) | Out-Null

# Force a run boundary before "explicitly instructs..." using a throwaway
# bookmark (bookmarks cannot live inside a run, so adding/removing one here
# splits the surrounding text into separate runs without altering content).
$preRange = $d.Content
$preRange.Find.Execute("explicitly instructs the character to do something", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$preRange.Collapse(1)
$d.Bookmarks.Add("ZZZTempSplit", $preRange)

# Place the real `_GoBack` bookmark right after "...do something"
$postRange = $d.Content
$postRange.Find.Execute("explicitly instructs the character to do something", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$postRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $postRange)

# Drop the scaffolding bookmark now that the run split exists
$d.Bookmarks("ZZZTempSplit").Delete()
